# Add ODK version form in exports
#
# Append a new "FormatVersion" / "format_version" row (row 24) to the
# old/new/odk_ref dictionary table on Sheet1, matching columns A (old),
# B (new) and C (odk_ref), and move the active selection the way the
# author left it (G12).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row values: old="FormatVersion", new="format_version", odk_ref="FormatVersion"
$ws.Range("A24").Value = "FormatVersion"
$ws.Range("B24").Value = "format_version"
$ws.Range("C24").Value = "FormatVersion"

# Columns A and B already default to the centered "no fill" look used
# throughout the table; give column C (which otherwise defaults to a plain
# style) the same centered formatting by copying it over from A24.
$ws.Range("A24").Copy()
$ws.Range("C24").PasteSpecial(-4122)

# Move the selection, as left by the author after the edit.
$null = $ws.Range("G12").Select()
